# Update Table 2 joinpoint estimates with more precise 2016 data
# (AAPC / APC / 95% CI values gain a third decimal place).
$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(3, 2).Range.Text = "0.741 (0.705, 0.776)"
$tbl.Cell(3, 3).Range.Text = "4.826 (4.736, 4.916)"
$tbl.Cell(3, 4).Range.Text = "13.042 (10.238, 15.917)*"
$tbl.Cell(3, 6).Range.Text = "0.765 (-3.198, 4.889)"
$tbl.Cell(3, 8).Range.Text = "11.521 (4.555, 18.950)*"
$tbl.Cell(3, 10).Range.Text = "33.471 (24.796, 42.750)*"
$tbl.Cell(4, 2).Range.Text = "0.833 (0.733, 0.933)"
$tbl.Cell(4, 3).Range.Text = "3.115 (2.946, 3.284)"
$tbl.Cell(4, 4).Range.Text = "10.435 (8.574, 12.327)*"
$tbl.Cell(4, 6).Range.Text = "-0.671 (-2.876, 1.585)"
$tbl.Cell(4, 8).Range.Text = "34.119 (29.704, 38.684)*"
$tbl.Cell(6, 2).Range.Text = "0.313 (0.290, 0.336)"
$tbl.Cell(6, 3).Range.Text = "1.206 (1.161, 1.250)"
$tbl.Cell(6, 4).Range.Text = "8.339 (5.706, 11.037)*"
$tbl.Cell(6, 6).Range.Text = "42.818 (33.145, 53.194)*"
$tbl.Cell(6, 8).Range.Text = "21.057 (6.972, 36.997)*"
$tbl.Cell(6, 10).Range.Text = "-6.174 (-7.091, -5.248)*"
$tbl.Cell(7, 2).Range.Text = "0.213 (0.163, 0.263)"
$tbl.Cell(7, 3).Range.Text = "0.579 (0.507, 0.651)"
$tbl.Cell(7, 4).Range.Text = "6.894 (2.691, 11.269)*"
$tbl.Cell(7, 6).Range.Text = "23.745 (18.357, 29.379)*"
$tbl.Cell(7, 8).Range.Text = "-10.324 (-28.427, 12.359)"
$tbl.Cell(7, 10).Range.Text = "-0.443 (-3.456, 2.663)"
$tbl.Cell(9, 2).Range.Text = "1.109 (1.066, 1.153)"
$tbl.Cell(9, 3).Range.Text = "4.609 (4.522, 4.695)"
$tbl.Cell(9, 4).Range.Text = "9.101 (7.958, 10.256)*"
$tbl.Cell(9, 6).Range.Text = "13.405 (11.647, 15.192)*"
$tbl.Cell(9, 8).Range.Text = "3.233 (1.471, 5.026)*"
$tbl.Cell(10, 2).Range.Text = "0.622 (0.536, 0.709)"
$tbl.Cell(10, 3).Range.Text = "2.113 (1.972, 2.254)"
$tbl.Cell(10, 4).Range.Text = "9.203 (8.152, 10.264)*"
$tbl.Cell(10, 6).Range.Text = "7.931 (7.326, 8.540)*"
$tbl.Cell(10, 8).Range.Text = "19.235 (10.047, 29.190)*"
$tbl.Cell(12, 2).Range.Text = "0.298 (0.275, 0.320)"
$tbl.Cell(12, 3).Range.Text = "3.591 (3.514, 3.669)"
$tbl.Cell(12, 4).Range.Text = "20.585 (17.159, 24.112)*"
$tbl.Cell(12, 6).Range.Text = "17.882 (11.531, 24.594)*"
$tbl.Cell(12, 8).Range.Text = "2.750 (-2.042, 7.776)"
$tbl.Cell(12, 10).Range.Text = "84.701 (71.656, 98.738)*"
$tbl.Cell(13, 2).Range.Text = "0.115 (0.078, 0.152)"
$tbl.Cell(13, 3).Range.Text = "2.102 (1.963, 2.240)"
$tbl.Cell(13, 4).Range.Text = "20.519 (8.923, 33.348)*"
$tbl.Cell(13, 6).Range.Text = "5.108 (-6.160, 17.730)"
$tbl.Cell(13, 8).Range.Text = "128.209 (61.589, 222.296)*"
$tbl.Cell(15, 2).Range.Text = "0.964 (0.923, 1.004)"
$tbl.Cell(15, 3).Range.Text = "1.007 (0.966, 1.047)"
$tbl.Cell(15, 4).Range.Text = "-0.237 (-0.752, 0.280)"
$tbl.Cell(15, 6).Range.Text = "-0.237 (-0.752, 0.280)"
$tbl.Cell(16, 2).Range.Text = "2.056 (1.900, 2.212)"
$tbl.Cell(16, 3).Range.Text = "0.594 (0.519, 0.668)"
$tbl.Cell(16, 4).Range.Text = "-5.796 (-9.374, -2.078)*"
$tbl.Cell(16, 6).Range.Text = "-8.476 (-9.589, -7.349)*"
$tbl.Cell(16, 8).Range.Text = "16.965 (-17.653, 66.137)"
